# Apply numeric recalculation updates to the Atomos_Profits leve-profit sheets
# (scheduled runner refresh of market-price-derived columns H:N)
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 37821.637
$ws.Range("J17").Value = 37821.637
$ws.Range("L17").Value = 113464.911
$ws.Range("N17").Value = -113800.911
$ws.Range("H33").Value = 283.53125
$ws.Range("I33").Value = 102.16
$ws.Range("K33").Value = 102.16
$ws.Range("M33").Value = 126.84
$ws.Range("H58").Value = 20881.963
$ws.Range("I58").Value = 258.22223
$ws.Range("J58").Value = 25100.455
$ws.Range("K58").Value = 774.66669
$ws.Range("L58").Value = 75301.36500000001
$ws.Range("M58").Value = -624.66669
$ws.Range("N58").Value = -75601.36500000001
$ws.Range("H69").Value = 3839
$ws.Range("I69").Value = 3863
$ws.Range("J69").Value = 3815
$ws.Range("K69").Value = 11589
$ws.Range("L69").Value = 11445
$ws.Range("M69").Value = -10715
$ws.Range("N69").Value = -13193
$ws.Range("H72").Value = 3839
$ws.Range("I72").Value = 3863
$ws.Range("J72").Value = 3815
$ws.Range("K72").Value = 34767
$ws.Range("L72").Value = 34335
$ws.Range("M72").Value = -30399
$ws.Range("N72").Value = -43071
$ws.Range("H98").Value = 2678.111
$ws.Range("I98").Value = 1557.6666
$ws.Range("J98").Value = 4919
$ws.Range("K98").Value = 1557.6666
$ws.Range("L98").Value = 4919
$ws.Range("M98").Value = -59.66660000000002
$ws.Range("N98").Value = -7915
$ws.Range("H122").Value = 2678.111
$ws.Range("I122").Value = 1557.6666
$ws.Range("J122").Value = 4919
$ws.Range("K122").Value = 4672.9998
$ws.Range("L122").Value = 14757
$ws.Range("M122").Value = -2222.9998
$ws.Range("N122").Value = -19657
$ws.Range("H129").Value = 13890114
$ws.Range("J129").Value = 1165.75
$ws.Range("L129").Value = 3497.25
$ws.Range("N129").Value = -13497.25
$ws.Range("H135").Value = 779.1778
$ws.Range("I135").Value = 549.4103
$ws.Range("J135").Value = 2272.6667
$ws.Range("K135").Value = 4944.6927
$ws.Range("L135").Value = 20454.0003
$ws.Range("M135").Value = -2409.6927
$ws.Range("N135").Value = -25524.0003
$ws.Range("H137").Value = 2396.8462
$ws.Range("I137").Value = 2504.8948
$ws.Range("J137").Value = 2103.5715
$ws.Range("K137").Value = 7514.6844
$ws.Range("L137").Value = 6310.7145
$ws.Range("M137").Value = -4964.6844
$ws.Range("N137").Value = -11410.7145
$ws.Range("H138").Value = 5201.3076
$ws.Range("I138").Value = 2405.9167
$ws.Range("J138").Value = 6039.925
$ws.Range("K138").Value = 7217.750100000001
$ws.Range("L138").Value = 18119.775
$ws.Range("M138").Value = -2077.750100000001
$ws.Range("N138").Value = -28399.775
$ws.Range("H141").Value = 606688.25
$ws.Range("I141").Value = 1689.4117
$ws.Range("J141").Value = 5749178.5
$ws.Range("K141").Value = 5068.2351
$ws.Range("L141").Value = 17247535.5
$ws.Range("M141").Value = 111.7649000000001
$ws.Range("N141").Value = -17257895.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 963.5417
$ws.Range("I74").Value = 902.6923
$ws.Range("K74").Value = 902.6923
$ws.Range("M74").Value = -28.69230000000005
$ws.Range("H77").Value = 963.5417
$ws.Range("I77").Value = 902.6923
$ws.Range("K77").Value = 4513.4615
$ws.Range("M77").Value = -145.4615000000003
$ws.Range("H132").Value = 23812060
$ws.Range("I132").Value = 43480190
$ws.Range("J132").Value = 3269.0527
$ws.Range("K132").Value = 130440570
$ws.Range("L132").Value = 9807.158100000001
$ws.Range("M132").Value = -130438040
$ws.Range("N132").Value = -14867.1581
$ws.Range("H137").Value = 30000
$ws.Range("J137").Value = 30000
$ws.Range("L137").Value = 30000
$ws.Range("N137").Value = -40200
$ws.Range("H141").Value = 30000
$ws.Range("J141").Value = 30000
$ws.Range("L141").Value = 30000
$ws.Range("N141").Value = -40360

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1592.2115
$ws.Range("I134").Value = 1197.766
$ws.Range("J134").Value = 5300
$ws.Range("K134").Value = 3593.298
$ws.Range("L134").Value = 15900
$ws.Range("M134").Value = -1058.298
$ws.Range("N134").Value = -20970

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2853.5745
$ws.Range("I31").Value = 1748.4546
$ws.Range("J31").Value = 3826.08
$ws.Range("K31").Value = 1748.4546
$ws.Range("L31").Value = 3826.08
$ws.Range("M31").Value = -1453.4546
$ws.Range("N31").Value = -4416.08
$ws.Range("H34").Value = 2853.5745
$ws.Range("I34").Value = 1748.4546
$ws.Range("J34").Value = 3826.08
$ws.Range("K34").Value = 1748.4546
$ws.Range("L34").Value = 3826.08
$ws.Range("M34").Value = -1546.4546
$ws.Range("N34").Value = -4230.08
$ws.Range("H58").Value = 9436356
$ws.Range("I58").Value = 1441.9459
$ws.Range("K58").Value = 1441.9459
$ws.Range("M58").Value = -1238.9459
$ws.Range("H136").Value = 9436356
$ws.Range("I136").Value = 1441.9459
$ws.Range("K136").Value = 4325.8377
$ws.Range("M136").Value = -1775.8377
$ws.Range("H137").Value = 29966.666
$ws.Range("J137").Value = 29966.666
$ws.Range("L137").Value = 29966.666
$ws.Range("N137").Value = -40166.666

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 11640.4
$ws.Range("I87").Value = 5251
$ws.Range("K87").Value = 15753
$ws.Range("M87").Value = -14505
$ws.Range("H90").Value = 11640.4
$ws.Range("I90").Value = 5251
$ws.Range("K90").Value = 47259
$ws.Range("M90").Value = -41019
$ws.Range("H118").Value = 1433.5834
$ws.Range("I118").Value = 401.14285
$ws.Range("J118").Value = 2879
$ws.Range("K118").Value = 1203.42855
$ws.Range("L118").Value = 8637
$ws.Range("M118").Value = 39.57144999999991
$ws.Range("N118").Value = -11123
$ws.Range("H120").Value = 13655.714
$ws.Range("I120").Value = 11118
$ws.Range("K120").Value = 33354
$ws.Range("M120").Value = -28516
$ws.Range("H121").Value = 28765.455
$ws.Range("I121").Value = 209.5
$ws.Range("J121").Value = 35111.223
$ws.Range("K121").Value = 628.5
$ws.Range("L121").Value = 105333.669
$ws.Range("M121").Value = 681.5
$ws.Range("N121").Value = -107953.669
$ws.Range("H125").Value = 2975
$ws.Range("J125").Value = 2975
$ws.Range("L125").Value = 8925
$ws.Range("N125").Value = -18765

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3427.75
$ws.Range("I113").Value = 1355.5
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 1355.5
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = 814.5
$ws.Range("N113").Value = -9840

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2819
$ws.Range("I132").Value = 1782.9445
$ws.Range("K132").Value = 5348.833500000001
$ws.Range("M132").Value = -2818.833500000001
$ws.Range("H136").Value = 1608.579
$ws.Range("I136").Value = 1005.25
$ws.Range("J136").Value = 2642.8572
$ws.Range("K136").Value = 3015.75
$ws.Range("L136").Value = 7928.571599999999
$ws.Range("M136").Value = -465.75
$ws.Range("N136").Value = -13028.5716

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4966.8335
$ws.Range("I62").Value = 4900.5
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 4900.5
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -4276.5
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 4966.8335
$ws.Range("I65").Value = 4900.5
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 24502.5
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -21382.5
$ws.Range("N65").Value = -31240
$ws.Range("H125").Value = 31687.5
$ws.Range("J125").Value = 31687.5
$ws.Range("L125").Value = 31687.5
$ws.Range("N125").Value = -41527.5
$ws.Range("H136").Value = 1071.0454
$ws.Range("I136").Value = 678.0968
$ws.Range("J136").Value = 2008.0769
$ws.Range("K136").Value = 2034.2904
$ws.Range("L136").Value = 6024.2307
$ws.Range("M136").Value = 515.7095999999999
$ws.Range("N136").Value = -11124.2307
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
